$d = $word.ActiveDocument

# The paragraph currently reads as two runs: "{m" + ":userdoc 'zone1'}".
# The target splits it into four runs: "{", "m", ":userdoc 'zone1'", "}".
# We locate the paragraph by its text (robust to any offset drift) and then
# force run boundaries at the three internal split points by temporarily
# bookmarking each point and immediately deleting the bookmark: adding a
# bookmark forces the underlying run to split at that character position,
# and removing the bookmark again leaves the split in place without adding
# any bookmark markers to the saved document.

function Split-RunAt($pos) {
    $r = $d.Range($pos, $pos)
    $bkName = "tmpSplitMarker"
    $d.Bookmarks.Add($bkName, $r) | Out-Null
    $d.Bookmarks($bkName).Delete()
}

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "{m:userdoc 'zone1'}`r") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph with text {m:userdoc 'zone1'}"
}

$start = $target.Range.Start

# Split points, expressed as offsets from the start of the paragraph range:
#   start+0 -> "{"
#   start+1 -> "m"
#   start+2 -> ":userdoc 'zone1'"
#   start+18 -> "}"
Split-RunAt ($start + 1)
Split-RunAt ($start + 2)
Split-RunAt ($start + 18)
